$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rating")

# All the player scores in column B are normalized to a single value (1600)
for ($row = 1; $row -le 18; $row++) {
    $ws.Cells.Item($row, 2).Value = 1600
}

# The computed rating values in column C are updated to reflect the new scores
$ws.Range("C1").Value = 1616
$ws.Range("C2").Value = 1584

# Move the active selection on this sheet
$ws.Activate()
$ws.Range("H17").Select()
